# Applies the updated cryptocurrency price/volume snapshot to Sheet1.
# Each target cell is explicitly formatted as Text ("@") before the
# value is written so strings such as "1.000" or "27.225.96" are not
# auto-coerced into numbers/dates by Excel, then the style is reset
# back to "Normal" so no residual cell formatting is introduced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cellUpdates = [ordered]@{
    "D2" = "27.225.96"
    "E2" = "  +0.21%  "
    "D3" = "1.905.65"
    "E3" = "  +0.71%  "
    "D4" = "1.000"
    "E4" = "  -0.28%  "
    "D5" = "306.22"
    "E5" = "  -0.26%  "
    "D6" = "1.000"
    "E6" = "  -0.21%  "
    "D7" = "0.5376"
    "E7" = "  +3.08%  "
    "D8" = "0.3807"
    "E8" = "  +1.44%  "
    "D9" = "0.07292"
    "E9" = "  +0.44%  "
    "D10" = "22.23"
    "E10" = "  +5.12%  "
    "D11" = "0.9052"
    "E11" = "  +0.75%  "
    "D12" = "0.08199"
    "E12" = "  +0.01%  "
    "D13" = "95.78"
    "E13" = "  -0.91%  "
    "D14" = "5.341"
    "E14" = "  +1.35%  "
    "D15" = "0.9991"
    "E15" = "  -0.35%  "
    "D16" = "14.86"
    "E16" = "  +2.26%  "
    "D17" = "0.000008667"
    "E17" = "  +0.92%  "
    "D18" = "1.000"
    "E18" = "  -0.17%  "
    "D19" = "27.263.39"
    "E19" = "  +0.27%  "
    "B20" = "WrappedEther"
    "C20" = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
    "D20" = "1.141.98"
    "E20" = "  -39.64%  "
    "B21" = "Uniswap"
    "C21" = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
    "D21" = "5.047"
    "E21" = "  -0.72%  "
    "E22" = "  +0.82%  "
    "D23" = "6.522"
    "E23" = "  +1.90%  "
    "D24" = "148.67"
    "E24" = "  +0.72%  "
    "E25" = "  +0.73%  "
    "E26" = "  +1.14%  "
    "D27" = "1.747"
    "D28" = "116.70"
    "D29" = "4.846"
    "E29" = "  +1.16%  "
    "D30" = "4.727"
    "E30" = "  -3.62%  "
    "D31" = "0.09220"
    "E31" = "  -0.08%  "
    "E32" = "  +4.83%  "
    "D33" = "0.05080"
    "E33" = "  +0.73%  "
    "E34" = "  +0.40%  "
    "D35" = "2.997"
    "E35" = "  +0.86%  "
    "D36" = "3.316"
    "E36" = "  -3.50%  "
    "E37" = "  +3.81%  "
    "D38" = "0.5955"
    "E38" = "  +5.18%  "
    "D39" = "0.02003"
    "E39" = "  +0.91%  "
    "E40" = "  +0.42%  "
    "D41" = "9.325"
    "E41" = "  +4.06%  "
    "D42" = "6.668"
    "E42" = "  +1.87%  "
    "D43" = "116.55"
    "E43" = "  +0.94%  "
    "D44" = "0.5163"
    "E44" = "  +6.32%  "
    "E45" = "  +0.83%  "
    "D46" = "10.22"
    "E46" = "  +1.16%  "
    "E47" = "  -0.22%  "
    "D48" = "1.644"
    "E48" = "  +1.43%  "
    "D49" = "38.37"
    "E49" = "  +0.70%  "
    "D50" = "0.06132"
    "E50" = "  +3.26%  "
    "D51" = "63.48"
}

foreach ($addr in $cellUpdates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $cellUpdates[$addr]
    $cell.Style = "Normal"
}
